# Add a new row (row 19) to Sheet1 for LeetCode problem 81:
# "Search in Rotated Sorted Array II"
# (mirrors the formatting of the row above it, row 18)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 19

# Copy formatting (styles / row height) from the last existing data row (18)
# onto the new row before filling in values.
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item($newRow).RowHeight = 34

# Fill in the new record's data
$ws.Cells.Item($newRow, 1).Value = 81
$ws.Cells.Item($newRow, 2).Value = "Search in Rotated Sorted Array II"
$ws.Cells.Item($newRow, 3).Value = "#array  #binary-search #必背 "
$ws.Cells.Item($newRow, 4).Value = "medium"
$ws.Cells.Item($newRow, 5).Value = 2
$ws.Cells.Item($newRow, 6).Value = 2
$ws.Cells.Item($newRow, 7).Value = 22
$ws.Cells.Item($newRow, 8).Value = (Get-Date -Year 2025 -Month 6 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($newRow, 9).Value = (Get-Date -Year 2025 -Month 6 -Day 27 -Hour 0 -Minute 0 -Second 0)

# Update the view: scroll position and active cell selection, matching
# where the user ended up after adding the new row.
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("I19").Select()
